$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Tempo (ms)" column B (rows 2-11)
$ws.Range("B2").Value = 94171.17857933044
$ws.Range("B3").Value = 79705.14178276062
$ws.Range("B4").Value = 75130.58233261108
$ws.Range("B5").Value = 64536.26704216003
$ws.Range("B6").Value = 72227.01478004456
$ws.Range("B7").Value = 69706.04991912842
$ws.Range("B8").Value = 67728.42264175415
$ws.Range("B9").Value = 49856.20522499084
$ws.Range("B10").Value = 49969.58160400391
$ws.Range("B11").Value = 45510.15019416809

# Update "Memória (KB)" column C (rows 2-11)
$ws.Range("C2").Value = -56
$ws.Range("C3").Value = -40
$ws.Range("C5").Value = 72
$ws.Range("C7").Value = -36
$ws.Range("C9").Value = 0
$ws.Range("C10").Value = 28
$ws.Range("C11").Value = -24

# Update statistics rows (Média / Mediana)
$ws.Range("B14").Value = 66854.05941009521
$ws.Range("C14").Value = -5.6
$ws.Range("B15").Value = 68717.23628044128
